# ex9.1.2(Linear) - "nuevos experimentos no convexos"
# Update the follower-restriction coefficients, the modified point, and the
# bf/BF vectors with the values of a new (non-convex) generated experiment.
#
# Sheet order in this workbook (1-based, matches xl/workbook.xml):
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha
# NOTE: worksheet name lookup via Worksheets.Item(name) is case-insensitive,
# and "Vector_bf" / "Vector_BF" would collide -- so every sheet below is
# addressed by its numeric index to stay unambiguous.

$wb = $excel.ActiveWorkbook

$sRestrFollower = $wb.Worksheets.Item(3)
$sPuntoMod      = $wb.Worksheets.Item(4)
$sVecbf         = $wb.Worksheets.Item(5)
$sVecBF         = $wb.Worksheets.Item(6)

# Cells whose new text looks like a plain number need to be forced to Text
# format first, otherwise Excel would silently store them as numeric values
# instead of keeping them as the literal text that belongs in this sheet.
$numericLookingCells = @(
    @($sRestrFollower, "B2"),
    @($sRestrFollower, "D2"),
    @($sRestrFollower, "F2"),
    @($sRestrFollower, "B3"),
    @($sRestrFollower, "D3"),
    @($sRestrFollower, "E3"),
    @($sRestrFollower, "F3"),
    @($sRestrFollower, "B4"),
    @($sRestrFollower, "D4"),
    @($sRestrFollower, "E4"),
    @($sRestrFollower, "F4"),
    @($sRestrFollower, "B5"),
    @($sRestrFollower, "D5"),
    @($sRestrFollower, "F5"),
    @($sPuntoMod,      "A2"),
    @($sPuntoMod,      "B2"),
    @($sVecbf,         "A2"),
    @($sVecBF,         "A2"),
    @($sVecBF,         "A3")
)

foreach ($pair in $numericLookingCells) {
    $pair[0].Range($pair[1]).NumberFormat = "@"
}

# Restricciones_del_follower ------------------------------------------------
$sRestrFollower.Range("A2").Value = "8.95 - y"
$sRestrFollower.Range("B2").Value = "-8.95"
$sRestrFollower.Range("D2").Value = "0.68"
$sRestrFollower.Range("F2").Value = "1.0"

$sRestrFollower.Range("A3").Value = "-1.9499999999999993 - x + y"
$sRestrFollower.Range("B3").Value = "-1.0500000000000007"
$sRestrFollower.Range("D3").Value = "0.24"
$sRestrFollower.Range("E3").Value = "0"
$sRestrFollower.Range("F3").Value = "6.4"

$sRestrFollower.Range("A4").Value = "-24.9 + x + 2y"
$sRestrFollower.Range("B4").Value = "12.899999999999999"
$sRestrFollower.Range("D4").Value = "0.44"
$sRestrFollower.Range("E4").Value = "0"
$sRestrFollower.Range("F4").Value = "5.2"

$sRestrFollower.Range("A5").Value = "-20.01 + 4x - y"
$sRestrFollower.Range("B5").Value = "7.050000000000001"
$sRestrFollower.Range("D5").Value = "1.0"
$sRestrFollower.Range("F5").Value = "5.8"

# Punto_modificado ------------------------------------------------------------
$sPuntoMod.Range("A2").Value = "7.0"
$sPuntoMod.Range("B2").Value = "8.95"

# Vector_bf ---------------------------------------------------------------
$sVecbf.Range("A2").Value = "-0.43999999999999995"

# Vector_BF ---------------------------------------------------------------
$sVecBF.Range("A2").Value = "1.0"
$sVecBF.Range("A3").Value = "3.0"
